# TestCasesForrester.xlsx — "stevan App IOS Automation"
#
# The "Navigate through onboarding screens" test case (Test Cases!D2) is
# updated from a "Pass" result to "Skipped", and the view state left behind
# by the editor (zoom level reset to 100%, and the last clicked cell on each
# sheet) is reproduced on both sheets.

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestData  = $wb.Worksheets.Item("Test Data")

# --- Content change -------------------------------------------------------
# Results column, "Navigate" test case: Pass -> Skipped
$wsTestCases.Range("D2").Value = "Skipped"

# --- View-state updates ----------------------------------------------------
# "Test Cases" sheet: zoom back to 100%, selection moved to H11
$wsTestCases.Activate()
$excel.ActiveWindow.Zoom = 100
$wsTestCases.Range("H11").Select()

# "Test Data" sheet: zoom back to 100%, selection moved to H15
# (this sheet remains the active/selected tab, as in the source file)
$wsTestData.Activate()
$excel.ActiveWindow.Zoom = 100
$wsTestData.Range("H15").Select()
